# Modify the contract of services: add LiveName / Result / ErrorMessage / LiveId
# fields to the StartLiveRequest / StartLiveResponse metadata table on sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlContinuous = 1
$xlLineStyleNone = -4142
$xlCenter = -4108
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

# ---------------------------------------------------------------------------
# 1. Grow the metadata table from 5 rows (StartLiveRequest: 2 fields,
#    StartLiveResponse: 1 field) to 7 rows (StartLiveRequest: 3 fields,
#    StartLiveResponse: 3 fields) by inserting two new blank rows.
#    Row 3 becomes a brand-new row (for LiveName) and, once that is in place,
#    row 6 becomes a second brand-new row (for LiveId).
# ---------------------------------------------------------------------------
$ws.Range("A3:E3").EntireRow.Insert()
$ws.Range("A6:E6").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Re-apply the merged cells for the (now taller) table.
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").Merge()
$ws.Range("B2:B4").Merge()
$ws.Range("B5:B7").Merge()
$ws.Range("A2:A7").Merge()

# ---------------------------------------------------------------------------
# Helper: apply the centred alignment + a box border (all four edges, or a
# chosen subset) used throughout this little table.
# ---------------------------------------------------------------------------
function Format-Cell($range, [bool]$left, [bool]$top, [bool]$bottom, [bool]$right) {
    if ($left)   { $range.Borders.Item($xlEdgeLeft).LineStyle   = $xlContinuous } else { $range.Borders.Item($xlEdgeLeft).LineStyle   = $xlLineStyleNone }
    if ($top)    { $range.Borders.Item($xlEdgeTop).LineStyle    = $xlContinuous } else { $range.Borders.Item($xlEdgeTop).LineStyle    = $xlLineStyleNone }
    if ($bottom) { $range.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous } else { $range.Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone }
    if ($right)  { $range.Borders.Item($xlEdgeRight).LineStyle  = $xlContinuous } else { $range.Borders.Item($xlEdgeRight).LineStyle  = $xlLineStyleNone }
    $range.HorizontalAlignment = $xlCenter
    $range.VerticalAlignment = $xlCenter
}

function Set-FullBox($range, $value) {
    if ($null -ne $value) { $range.Value = $value }
    Format-Cell $range $true $true $true $true
}

function Clear-FullBox($range) {
    $range.Value = ""
    Format-Cell $range $true $true $true $true
}

# ---------------------------------------------------------------------------
# Row 1 - header (B1:C1 is merged and drawn as a single box split into a
# left half (B1, no right edge) and a right half (C1, no left edge) so the
# merge doesn't show a divider line in the middle).
# ---------------------------------------------------------------------------
Set-FullBox $ws.Range("A1") "ServiceCode"
$ws.Range("B1").Value = "Name"
Format-Cell $ws.Range("B1") $true $true $true $false
$ws.Range("C1").Value = ""
Format-Cell $ws.Range("C1") $false $true $true $true
Set-FullBox $ws.Range("D1") "Metadata"
Set-FullBox $ws.Range("E1") "Desc"

# ---------------------------------------------------------------------------
# Row 2 - ServiceCode value + first StartLiveRequest field (LocalMac)
# ---------------------------------------------------------------------------
Set-FullBox $ws.Range("A2") 10000001
Set-FullBox $ws.Range("B2") "StartLiveRequest"
Set-FullBox $ws.Range("C2") "LocalMac"
Set-FullBox $ws.Range("D2") "String"
Set-FullBox $ws.Range("E2") "本机MAC"

# ---------------------------------------------------------------------------
# Row 3 (new) - StartLiveRequest field: LiveName
# ---------------------------------------------------------------------------
Clear-FullBox $ws.Range("A3")
Clear-FullBox $ws.Range("B3")
Set-FullBox $ws.Range("C3") "LiveName"
Set-FullBox $ws.Range("D3") "String"
Set-FullBox $ws.Range("E3") "直播的名称"

# ---------------------------------------------------------------------------
# Row 4 - StartLiveRequest field: FileSignature (this particular cell keeps
# the border-less look it already had before the edit)
# ---------------------------------------------------------------------------
Clear-FullBox $ws.Range("A4")
Clear-FullBox $ws.Range("B4")
$ws.Range("C4").Value = "FileSignature"
Format-Cell $ws.Range("C4") $false $false $false $false
Set-FullBox $ws.Range("D4") "String"
Set-FullBox $ws.Range("E4") "分享文件签名"

# ---------------------------------------------------------------------------
# Row 5 - StartLiveResponse value + first field: Result
# ---------------------------------------------------------------------------
Clear-FullBox $ws.Range("A5")
Set-FullBox $ws.Range("B5") "StartLiveResponse"
Set-FullBox $ws.Range("C5") "Result"
Set-FullBox $ws.Range("D5") "int"
Set-FullBox $ws.Range("E5") "服务结果：0成功；1失败"

# ---------------------------------------------------------------------------
# Row 6 (new) - StartLiveResponse field: LiveId
# ---------------------------------------------------------------------------
Clear-FullBox $ws.Range("A6")
Clear-FullBox $ws.Range("B6")
Set-FullBox $ws.Range("C6") "LiveId"
Set-FullBox $ws.Range("D6") "String"
Set-FullBox $ws.Range("E6") "直播签名，用户服务验证"

# ---------------------------------------------------------------------------
# Row 7 - StartLiveResponse field: ErrorMessage
# ---------------------------------------------------------------------------
Clear-FullBox $ws.Range("A7")
Clear-FullBox $ws.Range("B7")
Set-FullBox $ws.Range("C7") "ErrorMessage"
Set-FullBox $ws.Range("D7") "String"
Set-FullBox $ws.Range("E7") "失败原因描述"

# ---------------------------------------------------------------------------
# 3. Update the selection so it matches the saved workbook (cell C4).
# ---------------------------------------------------------------------------
$ws.Range("C4").Select()
